$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Remove the stray "rober " entry from cell A5 (leftover data from a
# previous search/reload). This also drops the now-unused shared string.
$ws.Range("A5").ClearContents()
